$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.749.03'
$ws.Range("E2").Value = '  -0.02%  '

$ws.Range("D3").Value = '3.793.48'
$ws.Range("E3").Value = '  -1.56%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.83'
$ws.Range("E5").Value = '  -0.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.83'
$ws.Range("E6").Value = '  +1.37%  '

$ws.Range("D7").Value = '3.794.20'
$ws.Range("E7").Value = '  -1.61%  '

$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.523'
$ws.Range("E9").Value = '  -0.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.165'
$ws.Range("E10").Value = '  +0.43%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.45'
$ws.Range("E11").Value = '  +1.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("E12").Value = '  +0.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000278'
$ws.Range("E13").Value = '  +12.51%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.52'
$ws.Range("E14").Value = '  -0.87%  '

$ws.Range("D15").Value = '4.428.76'
$ws.Range("E15").Value = '  -1.00%  '

$ws.Range("D16").Value = '3.816.86'
$ws.Range("E16").Value = '  -1.35%  '

$ws.Range("D17").Value = '67.792.86'
$ws.Range("E17").Value = '  +0.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.10'
$ws.Range("E18").Value = '  -0.16%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.34'
$ws.Range("E19").Value = '  -0.03%  '

$ws.Range("E20").Value = '  +0.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.72'
$ws.Range("E21").Value = '  -2.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '464.48'
$ws.Range("E22").Value = '  -0.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.719'
$ws.Range("E23").Value = '  -1.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000153'
$ws.Range("E24").Value = '  -5.86%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.06'
$ws.Range("E25").Value = '  -0.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.24'
$ws.Range("E26").Value = '  -0.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.02'
$ws.Range("E27").Value = '  -0.88%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.18'
$ws.Range("E28").Value = '  +1.76%  '

$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.91'
$ws.Range("E30").Value = '  -1.13%  '

$ws.Range("D31").Value = '3.952.95'
$ws.Range("E31").Value = '  -1.13%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.66'
$ws.Range("E32").Value = '  -1.10%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.27'
$ws.Range("E33").Value = '  -2.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.59'
$ws.Range("E34").Value = '  -1.59%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.27'
$ws.Range("E35").Value = '  -0.60%  '

$ws.Range("D36").Value = '3.759.54'
$ws.Range("E36").Value = '  -1.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.85'
$ws.Range("E37").Value = '  +18.36%  '

$ws.Range("E38").Value = '  +1.37%  '

$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.01'
$ws.Range("E39").Value = '  -1.42%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.139'
$ws.Range("E40").Value = '  -0.59%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.89'
$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.314'
$ws.Range("E43").Value = '  +0.86%  '

$ws.Range("B45").Value = 'Cosmos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.71'
$ws.Range("E45").Value = '  +2.45%  '

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.96'
$ws.Range("E46").Value = '  -0.58%  '

$ws.Range("B47").Value = 'FLOKI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000289'
$ws.Range("E47").Value = '  +5.34%  '

$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.30'
$ws.Range("E48").Value = '  -2.06%  '

$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '404.27'
$ws.Range("E49").Value = '  -5.62%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.43'
$ws.Range("E50").Value = '  -1.62%  '

$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0355'
$ws.Range("E51").Value = '  +0.48%  '
